$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.73"
$ws.Range("E2").Value = "'-3.43%"
$ws.Range("G2").Value = "'18"
$ws.Range("D3").Value = "'40.88"
$ws.Range("E3").Value = "'-2.21%"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'5.054"
$ws.Range("E4").Value = "'-2.67%"
$ws.Range("G4").Value = "'18"
$ws.Range("D5").Value = "'0.07608"
$ws.Range("E5").Value = "'-6.16%"
$ws.Range("G5").Value = "'18"
$ws.Range("E6").Value = "'-3.13%"
$ws.Range("G6").Value = "'18"
$ws.Range("D7").Value = "'1.595"
$ws.Range("E7").Value = "'-9.73%"
$ws.Range("G7").Value = "'18"
$ws.Range("E8").Value = "'-2.96%"
$ws.Range("G8").Value = "'18"
$ws.Range("D9").Value = "'0.09768"
$ws.Range("E9").Value = "'-12.90%"
$ws.Range("G9").Value = "'18"
$ws.Range("D10").Value = "'0.1765"
$ws.Range("E10").Value = "'-4.57%"
$ws.Range("G10").Value = "'18"
$ws.Range("D11").Value = "'0.09174"
$ws.Range("E11").Value = "'-0.87%"
$ws.Range("G11").Value = "'18"
$ws.Range("D12").Value = "'0.04353"
$ws.Range("E12").Value = "'-5.06%"
$ws.Range("G12").Value = "'18"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("E13").Value = "'-0.14%"
$ws.Range("G13").Value = "'18"
$ws.Range("D14").Value = "'0.001252"
$ws.Range("E14").Value = "'-1.60%"
$ws.Range("G14").Value = "'18"
$ws.Range("D15").Value = "'0.005793"
$ws.Range("E15").Value = "'-1.17%"
$ws.Range("G15").Value = "'18"
$ws.Range("E16").Value = "'0.83%"
$ws.Range("G16").Value = "'18"
$ws.Range("D17").Value = "'2.437"
$ws.Range("E17").Value = "'-8.29%"
$ws.Range("G17").Value = "'18"
$ws.Range("D18").Value = "'0.3279"
$ws.Range("E18").Value = "'-2.09%"
$ws.Range("G18").Value = "'18"
$ws.Range("D19").Value = "'6.762"
$ws.Range("E19").Value = "'-8.15%"
$ws.Range("G19").Value = "'18"
$ws.Range("E20").Value = "'-2.28%"
$ws.Range("G20").Value = "'18"
$ws.Range("E21").Value = "'6.91%"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'0.04163"
$ws.Range("E22").Value = "'-0.35%"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'0.001217"
$ws.Range("E23").Value = "'-1.98%"
$ws.Range("G23").Value = "'18"
$ws.Range("D24").Value = "'0.003995"
$ws.Range("E24").Value = "'-6.42%"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.0001300"
$ws.Range("E25").Value = "'6.39%"
$ws.Range("G25").Value = "'18"
$ws.Range("E26").Value = "'0.91%"
$ws.Range("G26").Value = "'18"
$ws.Range("G27").Value = "'18"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("D38").Value = "'0.02434"
$ws.Range("E38").Value = "'-5.86%"
$ws.Range("G38").Value = "'18"
$ws.Range("D39").Value = "'0.05146"
$ws.Range("E39").Value = "'-6.03%"
$ws.Range("G39").Value = "'18"
$ws.Range("D40").Value = "'0.007845"
$ws.Range("E40").Value = "'-2.26%"
$ws.Range("G40").Value = "'18"
$ws.Range("E41").Value = "'-6.42%"
$ws.Range("G41").Value = "'18"
$ws.Range("D42").Value = "'0.007078"
$ws.Range("E42").Value = "'-3.33%"
$ws.Range("G42").Value = "'18"
$ws.Range("D43").Value = "'0.001949"
$ws.Range("E43").Value = "'-6.44%"
$ws.Range("G43").Value = "'18"
$ws.Range("D44").Value = "'0.008379"
$ws.Range("E44").Value = "'1.61%"
$ws.Range("G44").Value = "'18"
$ws.Range("D45").Value = "'0.3335"
$ws.Range("E45").Value = "'-3.37%"
$ws.Range("G45").Value = "'18"
$ws.Range("D46").Value = "'0.00006376"
$ws.Range("E46").Value = "'-5.27%"
$ws.Range("G46").Value = "'18"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("G47").Value = "'18"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.003003"
$ws.Range("E48").Value = "'-26.83%"
$ws.Range("G48").Value = "'18"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.006436"
$ws.Range("E49").Value = "'89.80%"
$ws.Range("G49").Value = "'18"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("G50").Value = "'18"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.03%"
$ws.Range("G51").Value = "'18"
